# Append a new "update" log entry (row 39) plus its trailing 0-9 index row
# (row 40) to Sheet1, mirroring the existing blocks already on the sheet
# (e.g. rows 32/33 "tagsusp2tl", rows 35/36 "kirimlapts").
#
# Row 39 spells, across C39:I39, the word "update", followed by a "|"
# separator, a run of digits (J39:P39), an "F" flag (Q39), another "|"
# separator (R39) and a trailing digit (S39).
# Row 40 is the usual 0-9 filler row (C40:L40).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 39 : "update|890806|F|1" -----------------------------------------
$ws.Range("C39").Value = "u"
$ws.Range("D39").Value = "p"
$ws.Range("E39").Value = "d"
$ws.Range("F39").Value = "a"
$ws.Range("G39").Value = "t"
$ws.Range("H39").Value = "e"
$ws.Range("I39").Value = "|"
$ws.Range("J39").Value = 8
$ws.Range("K39").Value = 9
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = 8
$ws.Range("N39").Value = 0
$ws.Range("O39").Value = 8
$ws.Range("P39").Value = 6
$ws.Range("Q39").Value = "F"
$ws.Range("R39").Value = "|"
$ws.Range("S39").Value = 1

# --- Row 40 : trailing 0-9 filler row --------------------------------------
$ws.Range("C40").Value = 0
$ws.Range("D40").Value = 1
$ws.Range("E40").Value = 2
$ws.Range("F40").Value = 3
$ws.Range("G40").Value = 4
$ws.Range("H40").Value = 5
$ws.Range("I40").Value = 6
$ws.Range("J40").Value = 7
$ws.Range("K40").Value = 8
$ws.Range("L40").Value = 9

# --- View state: move the cursor to where it ended up after the edit ------
$ws.Range("AH34").Select()
$excel.ActiveWindow.ScrollRow = 12
$excel.ActiveWindow.ScrollColumn = 1
